# Apply the Ledger functionality TC updates described in the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newText = "It should not be displayed and should be as per parameter."

# Rows 19-20 (Step 12 / Step 13) - column F expected-result text
$ws.Range("F19").Value = $newText
$ws.Range("F20").Value = $newText

# Rows 57-58 - column F expected-result text (row height should revert to default
# now that the text is shorter and no longer needs to wrap)
$ws.Range("F57").Value = $newText
$ws.Range("F58").Value = $newText
$ws.Rows("57:58").AutoFit()

# Update the active selection to match the saved view state
$ws.Range("F57:F58").Select()
